$d = $word.ActiveDocument

# --- Part 0: the "_GoBack" bookmark currently lives by itself in an
#     empty paragraph near the end of the document. It is being moved
#     to the "minutos" paragraph below, so remove it from its old spot
#     first (bookmark names must stay unique, and this also matches the
#     second hunk of the diff, which simply deletes the bookmark pair).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- Part 1: find the paragraph whose whole text is "minutos" (the
#     estimated-time value paragraph) and turn it into "50 minutos",
#     with the "_GoBack" bookmark sitting right between the two words,
#     and without the old w:proofErr gramStart/gramEnd wrapper.
$target = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text.Trim() -eq "minutos") {
        $target = $para
        break
    }
}

if ($target -ne $null) {
    $full = $target.Range

    $xml = '<?xml version="1.0" standalone="yes"?>' +
      '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
      '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
      '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
      '<w:body>' +
      '<w:p w14:paraId="71B1C5EB" w14:textId="4E912351" w:rsidR="00CD652E" w:rsidRDefault="00733C27" w:rsidP="00CD652E">' +
      '<w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="es-ES_tradnl"/></w:rPr></w:pPr>' +
      '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t xml:space="preserve">50 </w:t></w:r>' +
      '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
      '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t>minutos</w:t></w:r>' +
      '</w:p>' +
      '</w:body></w:document>' +
      '</pkg:xmlData></pkg:part></pkg:package>'

    [void]$full.InsertXML($xml)
}
